$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.69827809267175
$ws.Range("C2").Value = 8.598051137455585
$ws.Range("E2").Value = 24.38119412457513
$ws.Range("F2").Value = 39.2854062845096
$ws.Range("G2").Value = 3.603970106440488
$ws.Range("I2").Value = 17.29800784625574
$ws.Range("J2").Value = 7.529474323723931
$ws.Range("N2").Value = 15.96214754437588
$ws.Range("O2").Value = 18.65665522561837
$ws.Range("B3").Value = 12.0980240920947
$ws.Range("C3").Value = 8.081893876172835
$ws.Range("E3").Value = 24.13200284781451
$ws.Range("F3").Value = 39.03710725320976
$ws.Range("G3").Value = 3.605968780589484
$ws.Range("I3").Value = 17.40951418178398
$ws.Range("J3").Value = 7.553038055019439
$ws.Range("N3").Value = 15.9900872514857
$ws.Range("O3").Value = 18.71058814063014
$ws.Range("B4").Value = 11.71453958006393
$ws.Range("C4").Value = 7.74611239302687
$ws.Range("E4").Value = 23.98327575748047
$ws.Range("F4").Value = 38.89574261194738
$ws.Range("G4").Value = 3.607260174729513
$ws.Range("I4").Value = 17.48244855097002
$ws.Range("J4").Value = 7.568463883877579
$ws.Range("N4").Value = 16.00903487441506
$ws.Range("O4").Value = 18.74944699478321
$ws.Range("B5").Value = 11.55470569849868
$ws.Range("C5").Value = 7.604562852157555
$ws.Range("E5").Value = 23.92380453275681
$ws.Range("F5").Value = 38.84097356946033
$ws.Range("G5").Value = 3.607802624653147
$ws.Range("I5").Value = 17.51329168892316
$ws.Range("J5").Value = 7.574991110766714
$ws.Range("N5").Value = 16.01720762302181
$ws.Range("O5").Value = 18.76671954434211
$ws.Range("B6").Value = 11.52795655616836
$ws.Range("C6").Value = 7.580774773725822
$ws.Range("E6").Value = 23.91399973618723
$ws.Range("F6").Value = 38.83205196115323
$ws.Range("G6").Value = 3.607893677817174
$ws.Range("I6").Value = 17.51848086574191
$ws.Range("J6").Value = 7.576089520563952
$ws.Range("N6").Value = 16.01859198768265
$ws.Range("O6").Value = 18.76967422927428
$ws.Range("B7").Value = 11.71239814577406
$ws.Range("C7").Value = 7.744222459716068
$ws.Range("E7").Value = 23.98246902863541
$ws.Range("F7").Value = 38.89499242497307
$ws.Range("G7").Value = 3.607267424750783
$ws.Range("I7").Value = 17.48285997281032
$ws.Range("J7").Value = 7.568550935857757
$ws.Range("N7").Value = 16.00914326620062
$ws.Range("O7").Value = 18.74967412876469
$ws.Range("B8").Value = 12.49451488629151
$ws.Range("C8").Value = 8.424007516803002
$ws.Range("E8").Value = 24.29442231367648
$ws.Range("F8").Value = 39.19752164110001
$ws.Range("G8").Value = 3.604645956324924
$ws.Range("I8").Value = 17.3355264706983
$ws.Range("J8").Value = 7.537400515405812
$ws.Range("N8").Value = 15.97140953184898
$ws.Range("O8").Value = 18.67405503147092
$ws.Range("B9").Value = 13.90277160189843
$ws.Range("C9").Value = 9.60666440794593
$ws.Range("E9").Value = 24.93733178402937
$ws.Range("F9").Value = 39.87639125713677
$ws.Range("G9").Value = 3.60001225622133
$ws.Range("I9").Value = 17.08215660048987
$ws.Range("J9").Value = 7.483900214485847
$ws.Range("N9").Value = 15.91160668870446
$ws.Range("O9").Value = 18.57163865212575
$ws.Range("B10").Value = 14.85321527963026
$ws.Range("C10").Value = 10.38328970993891
$ws.Range("E10").Value = 25.42485896349632
$ws.Range("F10").Value = 40.42380954960214
$ws.Range("G10").Value = 3.596913576476826
$ws.Range("I10").Value = 16.91780673557703
$ws.Range("J10").Value = 7.449200598129911
$ws.Range("N10").Value = 15.87628206912536
$ws.Range("O10").Value = 18.52473659275564
$ws.Range("B11").Value = 15.26610023605101
$ws.Range("C11").Value = 10.71652972992137
$ws.Range("E11").Value = 25.64911146190335
$ws.Range("F11").Value = 40.68256916519324
$ws.Range("G11").Value = 3.595569567546763
$ws.Range("I11").Value = 16.84780356889368
$ws.Range("J11").Value = 7.434411677072434
$ws.Range("N11").Value = 15.86207360707076
$ws.Range("O11").Value = 18.50962451653743
$ws.Range("B12").Value = 15.41956949032859
$ws.Range("C12").Value = 10.83983669686404
$ws.Range("E12").Value = 25.73431549803811
$ws.Range("F12").Value = 40.7818790724148
$ws.Range("G12").Value = 3.595070004427329
$ws.Range("I12").Value = 16.82198256019504
$ws.Range("J12").Value = 7.428954492171572
$ws.Range("N12").Value = 15.85696011221095
$ws.Range("O12").Value = 18.50480182270367
$ws.Range("B13").Value = 15.38664642171898
$ws.Range("C13").Value = 10.8134085744469
$ws.Range("E13").Value = 25.71595376296618
$ws.Range("F13").Value = 40.76043330677344
$ws.Range("G13").Value = 3.59517717770435
$ws.Range("I13").Value = 16.82751294078387
$ws.Range("J13").Value = 7.430123434806152
$ws.Range("N13").Value = 15.85804953300146
$ws.Range("O13").Value = 18.5058003797494
$ws.Range("B14").Value = 15.27878440130595
$ws.Range("C14").Value = 10.72673205787837
$ws.Range("E14").Value = 25.65611595480864
$ws.Range("F14").Value = 40.69071338070314
$ws.Range("G14").Value = 3.595528280399366
$ws.Range("I14").Value = 16.84566545861037
$ws.Range("J14").Value = 7.433959845086132
$ws.Range("N14").Value = 15.86164757063328
$ws.Range("O14").Value = 18.50920968874361
$ws.Range("B15").Value = 15.21233841719592
$ws.Range("C15").Value = 10.67326466135602
$ws.Range("E15").Value = 25.61949850779066
$ws.Range("F15").Value = 40.64817791311333
$ws.Range("G15").Value = 3.595744561557396
$ws.Range("I15").Value = 16.85687405018354
$ws.Range("J15").Value = 7.43632838237767
$ws.Range("N15").Value = 15.86388621672685
$ws.Range("O15").Value = 18.51141532136985
$ws.Range("B16").Value = 14.82583313642706
$ws.Range("C16").Value = 10.3611081141432
$ws.Range("E16").Value = 25.41024721925471
$ws.Range("F16").Value = 40.40708841815931
$ws.Range("G16").Value = 3.597002726406669
$ws.Range("I16").Value = 16.92247766331446
$ws.Range("J16").Value = 7.450187117525518
$ws.Range("N16").Value = 15.87724801579274
$ws.Range("O16").Value = 18.52584987857374
$ws.Range("B17").Value = 14.58367421589368
$ws.Range("C17").Value = 10.16447348554256
$ws.Range("E17").Value = 25.28246119343406
$ws.Range("F17").Value = 40.26162779092885
$ws.Range("G17").Value = 3.597791335719684
$ws.Range("I17").Value = 16.96394479556932
$ws.Range("J17").Value = 7.458943998856045
$ws.Range("N17").Value = 15.88592120309744
$ws.Range("O17").Value = 18.53630280116307
$ws.Range("B18").Value = 14.44256139445571
$ws.Range("C18").Value = 10.04948702229495
$ws.Range("E18").Value = 25.20919917278948
$ws.Range("F18").Value = 40.17888390643878
$ws.Range("G18").Value = 3.598251100038169
$ws.Range("I18").Value = 16.98824347017229
$ws.Range("J18").Value = 7.464074494121673
$ws.Range("N18").Value = 15.89108498575515
$ws.Range("O18").Value = 18.54290061127469
$ws.Range("B19").Value = 14.39447150879365
$ws.Range("C19").Value = 10.01023047391509
$ws.Range("E19").Value = 25.184436812617
$ws.Range("F19").Value = 40.15102876283938
$ws.Range("G19").Value = 3.598407830750127
$ws.Range("I19").Value = 16.99654741401879
$ws.Range("J19").Value = 7.465827702172977
$ws.Range("N19").Value = 15.89286346342293
$ws.Range("O19").Value = 18.54523492787878
$ws.Range("B20").Value = 14.60964241332461
$ws.Range("C20").Value = 10.18560093086494
$ws.Range("E20").Value = 25.29604020993517
$ws.Range("F20").Value = 40.27701752745731
$ws.Range("G20").Value = 3.597706747959553
$ws.Range("I20").Value = 16.95948417268984
$ws.Range("J20").Value = 7.458002110510028
$ws.Range("N20").Value = 15.88497980018747
$ws.Range("O20").Value = 18.5351294282139
$ws.Range("B21").Value = 15.31054486172424
$ws.Range("C21").Value = 10.75226930887352
$ws.Range("E21").Value = 25.67368461509021
$ws.Range("F21").Value = 40.71115652594594
$ws.Range("G21").Value = 3.595424898737946
$ws.Range("I21").Value = 16.84031493529335
$ws.Range("J21").Value = 7.432829117223115
$ws.Range("N21").Value = 15.86058350028808
$ws.Range("O21").Value = 18.50818383212345
$ws.Range("B22").Value = 15.75180277702248
$ws.Range("C22").Value = 11.10581232598255
$ws.Range("E22").Value = 25.92212251003643
$ws.Range("F22").Value = 41.00256696012404
$ws.Range("G22").Value = 3.593988252495772
$ws.Range("I22").Value = 16.76644077155879
$ws.Range("J22").Value = 7.417210887448793
$ws.Range("N22").Value = 15.84619473773157
$ws.Range("O22").Value = 18.49582061580908
$ws.Range("B23").Value = 15.51785743414259
$ws.Range("C23").Value = 10.91865719353123
$ws.Range("E23").Value = 25.78940123950644
$ws.Range("F23").Value = 40.84635947664204
$ws.Range("G23").Value = 3.594750030959761
$ws.Range("I23").Value = 16.80550083532082
$ws.Range("J23").Value = 7.425470397607493
$ws.Range("N23").Value = 15.85373216851506
$ws.Range("O23").Value = 18.50193749742876
$ws.Range("B24").Value = 14.59790807949065
$ws.Range("C24").Value = 10.17605525113832
$ws.Range("E24").Value = 25.28990049607679
$ws.Range("F24").Value = 40.27005707001363
$ws.Range("G24").Value = 3.597744970201433
$ws.Range("I24").Value = 16.9614993917527
$ws.Range("J24").Value = 7.458427639036019
$ws.Range("N24").Value = 15.8854048556257
$ws.Range("O24").Value = 18.5356580779743
$ws.Range("B25").Value = 13.53611505895854
$ws.Range("C25").Value = 9.302905310005841
$ws.Range("E25").Value = 24.76046049315855
$ws.Range("F25").Value = 39.6839218618771
$ws.Range("G25").Value = 3.601211868195855
$ws.Range("I25").Value = 17.14688112501953
$ws.Range("J25").Value = 7.497563194865618
$ws.Range("N25").Value = 15.92626958708597
$ws.Range("O25").Value = 18.59439300510859
